# issue #5: stock data output to json file
#
# The "股票" (stock) worksheet gets a new "property_category" column
# inserted between the existing "total" and "date" columns (i.e. before
# the current column H), with value "stock" in every data row. This
# pushes the former H/I/J ("date" / "legislator_name" / "legislator_id")
# columns one slot to the right (I/J/K).
#
# A handful of company-name strings on that same sheet also get an
# accidental embedded space removed (OCR/typo cleanup).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)   # "股票" sheet

# Insert a new column before column H (the current "date" column),
# shifting date/legislator_name/legislator_id one column to the right.
$ws.Columns.Item(8).Insert()

# New header for the inserted column.
$ws.Range("H1").Value = "property_category"

# Fill the new column with the "stock" category marker for every data row.
$lastRow = 11
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 8).Value = "stock"
}

# Clean up the company names that had a stray embedded space.
$ws.Range("B4").Value = "南亞塑膠工業股份有限公司"
$ws.Range("B6").Value = "台新金融控股股份有限公司"
$ws.Range("B7").Value = "中國信託金融控股股份有限公司"
$ws.Range("B8").Value = "台灣苯乙烯工業股份有限公司"
$ws.Range("B9").Value = "第一金融控股股份有限公司"
$ws.Range("B11").Value = "新光合成纖维股份有限公司"
